# Update gh-pages to output generated at 456a3b4
#
# Sheet "展览" (sheet1): insert one new event row (2024-10-26, 亿万心动国乙)
# before the old row 36, shifting the following 3 rows down, and refresh a
# number of "想去人数" (F column) counters that simply ticked up since the
# last crawl. The same refreshed counters are mirrored onto the "演出",
# "本地生活" and "全部类型" sheets wherever the same event is duplicated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Insert a new row at position 36, pushing the old rows 36-38 down to 37-39.
$ws1.Rows.Item(36).Insert()

# Copy the (now shifted) A37 cell's formatting back onto the new A36 cell so
# it keeps the same bordered/centered/bold style used by every other row in
# column A.
$ws1.Cells.Item(37, 1).Copy()
$ws1.Cells.Item(36, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the brand new row describing the "亿万心动国乙" event.
$ws1.Cells.Item(36, 1).Value = 35
$ws1.Range("B36").Value = "'2024-10-26"
$ws1.Range("C36").Value = "杭州·亿万心动国乙✘代号鸢同人only(日夜场）"
$ws1.Range("D36").Value = "皓月路299号 诺丁山艺术中心"
$ws1.Range("E36").Value = "2024.10.26 10:00-10.26 21:00"
$ws1.Cells.Item(36, 6).Value = 4
$ws1.Cells.Item(36, 7).Value = 75
$ws1.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=91962"
$ws1.Range("I36").Value = "//i1.hdslb.com/bfs/openplatform/202409/dd2vRpMx1725458690845.jpeg"

# The "index" column (A) is just row-number minus one; Insert() shifted the
# rows down but left the old literal numbers behind, so the three rows that
# moved (old 36/37/38 -> new 37/38/39) need their index bumped by one.
$ws1.Cells.Item(37, 1).Value = 36
$ws1.Cells.Item(38, 1).Value = 37
$ws1.Cells.Item(39, 1).Value = 38

# Row 37 (was row 36, "BanGDream! Only同人展") also picked up more RSVPs.
$ws1.Cells.Item(37, 6).Value = 81

# Row 39 (was row 38, "原神X崩坏X星铁旅行盛宴") also picked up more RSVPs.
$ws1.Cells.Item(39, 6).Value = 3921

# Refresh the "想去人数" counters for events that did not move rows.
$ws1.Range("F2").Value = 1265
$ws1.Range("F5").Value = 5599
$ws1.Range("F6").Value = 1797
$ws1.Range("F7").Value = 1797
$ws1.Range("F8").Value = 6366
$ws1.Range("F10").Value = 1931
$ws1.Range("F11").Value = 516
$ws1.Range("F12").Value = 11
$ws1.Range("F18").Value = 7942
$ws1.Range("F19").Value = 7942
$ws1.Range("F24").Value = 1746
$ws1.Range("F26").Value = 3
$ws1.Range("F31").Value = 1767
$ws1.Range("F32").Value = 805
$ws1.Range("F33").Value = 376

# ---------------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 90

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 687
$ws3.Range("F5").Value = 271

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (aggregates every row from the three sheets above - no
# row is inserted here, only the mirrored counters are refreshed)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 687
$ws4.Range("F5").Value = 271
$ws4.Range("F10").Value = 5599
$ws4.Range("F11").Value = 271
$ws4.Range("F12").Value = 1797
$ws4.Range("F13").Value = 1797
$ws4.Range("F14").Value = 6366
$ws4.Range("F16").Value = 1931
$ws4.Range("F18").Value = 516
$ws4.Range("F24").Value = 7942
$ws4.Range("F25").Value = 7942
$ws4.Range("F30").Value = 1746
$ws4.Range("F32").Value = 3
$ws4.Range("F36").Value = 1767
$ws4.Range("F37").Value = 805
$ws4.Range("F39").Value = 376
$ws4.Range("F44").Value = 90
$ws4.Range("F47").Value = 3921
